# Auto-generated Word COM-interop script
# Replaces multiplication problem text in each table cell
$d = $word.ActiveDocument

function Replace-Exact($findText, $replaceText) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Replacement.ClearFormatting()
    $found = $range.Find.Execute($findText, $true, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $found) {
        Write-Host "WARNING: not found -> $findText"
    }
}

Replace-Exact '254×5=1270' '333×3=999'
Replace-Exact '921×2=1842' '708×8=5664'
Replace-Exact '933×4=3732' '999×3=2997'
Replace-Exact '435×9=3915' '119×8=952'
Replace-Exact '965×3=2895' '641×6=3846'
Replace-Exact '963×8=7704' '938×9=8442'
Replace-Exact '746×2=1492' '367×2=734'
Replace-Exact '917×5=4585' '909×4=3636'
Replace-Exact '289×8=2312' '110×2=220'
Replace-Exact '584×8=4672' '921×5=4605'
Replace-Exact '619×7=4333' '648×9=5832'
Replace-Exact '982×5=4910' '609×6=3654'
Replace-Exact '551×4=2204' '598×4=2392'
Replace-Exact '627×9=5643' '499×7=3493'
Replace-Exact '951×4=3804' '261×3=783'
Replace-Exact '441×8=3528' '819×6=4914'
Replace-Exact '436×8=3488' '282×9=2538'
Replace-Exact '511×3=1533' '488×6=2928'
Replace-Exact '699×2=1398' '589×4=2356'
Replace-Exact '992×8=7936' '960×2=1920'
Replace-Exact '441×7=3087' '911×9=8199'
Replace-Exact '563×4=2252' '907×4=3628'
Replace-Exact '291×8=2328' '642×5=3210'
Replace-Exact '365×5=1825' '666×5=3330'
Replace-Exact '133×2=266' '737×3=2211'

Write-Host "Done replacing $($d.Content.Text.Length) chars of content"
